$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32
$ws.Cells.Item(32, 1).Value = 131106660
$ws.Cells.Item(32, 2).Value = 79000
$ws.Cells.Item(32, 4).Value = "'NT"
$ws.Cells.Item(32, 5).Value = 6446
$ws.Cells.Item(32, 6).Value = "'Kolflarnlav"
$ws.Cells.Item(32, 7).Value = "'Carbonicola anthracophila"
$ws.Cells.Item(32, 8).Value = "'(Nyl.) Bendiksby & Timdal"
$ws.Cells.Item(32, 9).Value = "'1"
$ws.Cells.Item(32, 10).Value = "'dm²"
$ws.Cells.Item(32, 16).Value = "'Lappberget, Mpd"
$ws.Cells.Item(32, 17).Value = 601244
$ws.Cells.Item(32, 18).Value = 6959831
$ws.Cells.Item(32, 19).Value = 10
$ws.Cells.Item(32, 20).Value = "'Västernorrland"
$ws.Cells.Item(32, 21).Value = "'Timrå"
$ws.Cells.Item(32, 22).Value = "'Medelpad"
$ws.Cells.Item(32, 23).Value = "'Ljustorp"
$ws.Cells.Item(32, 24).Value = "'2025_0516"
$ws.Cells.Item(32, 25).Value = "'2025-06-26"
$ws.Cells.Item(32, 26).Value = "'08:28"
$ws.Cells.Item(32, 27).Value = "'2025-06-26"
$ws.Cells.Item(32, 28).Value = "'08:28"
$ws.Cells.Item(32, 30).Value = $false
$ws.Cells.Item(32, 31).Value = $false
$ws.Cells.Item(32, 33).Value = $false
$ws.Cells.Item(32, 46).Value = "'"
$ws.Cells.Item(32, 49).Value = "'David Isaksson"
$ws.Cells.Item(32, 50).Value = "'Samuel Koont"
$ws.Cells.Item(32, 51).Value = "'Kustpaketet"

# Row 33
$ws.Cells.Item(33, 1).Value = 131106663
$ws.Cells.Item(33, 2).Value = 79862
$ws.Cells.Item(33, 4).Value = "'NT"
$ws.Cells.Item(33, 5).Value = 6453
$ws.Cells.Item(33, 6).Value = "'Vedskivlav"
$ws.Cells.Item(33, 7).Value = "'Hertelidea botryosa"
$ws.Cells.Item(33, 8).Value = "'(Fr.) Printzen & Kantvilas"
$ws.Cells.Item(33, 9).Value = "'1"
$ws.Cells.Item(33, 10).Value = "'dm²"
$ws.Cells.Item(33, 16).Value = "'Lappberget, Mpd"
$ws.Cells.Item(33, 17).Value = 601228
$ws.Cells.Item(33, 18).Value = 6959812
$ws.Cells.Item(33, 19).Value = 10
$ws.Cells.Item(33, 20).Value = "'Västernorrland"
$ws.Cells.Item(33, 21).Value = "'Timrå"
$ws.Cells.Item(33, 22).Value = "'Medelpad"
$ws.Cells.Item(33, 23).Value = "'Ljustorp"
$ws.Cells.Item(33, 24).Value = "'2025_0513"
$ws.Cells.Item(33, 25).Value = "'2025-06-26"
$ws.Cells.Item(33, 26).Value = "'08:22"
$ws.Cells.Item(33, 27).Value = "'2025-06-26"
$ws.Cells.Item(33, 28).Value = "'08:22"
$ws.Cells.Item(33, 29).Value = "'tallstubbe"
$ws.Cells.Item(33, 30).Value = $false
$ws.Cells.Item(33, 31).Value = $false
$ws.Cells.Item(33, 33).Value = $false
$ws.Cells.Item(33, 46).Value = "'"
$ws.Cells.Item(33, 49).Value = "'David Isaksson"
$ws.Cells.Item(33, 50).Value = "'Måns Svensson"
$ws.Cells.Item(33, 51).Value = "'Kustpaketet"

# Row 34
$ws.Cells.Item(34, 1).Value = 131106648
$ws.Cells.Item(34, 2).Value = 79862
$ws.Cells.Item(34, 4).Value = "'NT"
$ws.Cells.Item(34, 5).Value = 6453
$ws.Cells.Item(34, 6).Value = "'Vedskivlav"
$ws.Cells.Item(34, 7).Value = "'Hertelidea botryosa"
$ws.Cells.Item(34, 8).Value = "'(Fr.) Printzen & Kantvilas"
$ws.Cells.Item(34, 9).Value = "'"
$ws.Cells.Item(34, 16).Value = "'Lappberget, Mpd"
$ws.Cells.Item(34, 17).Value = 601198
$ws.Cells.Item(34, 18).Value = 6959731
$ws.Cells.Item(34, 19).Value = 10
$ws.Cells.Item(34, 20).Value = "'Västernorrland"
$ws.Cells.Item(34, 21).Value = "'Timrå"
$ws.Cells.Item(34, 22).Value = "'Medelpad"
$ws.Cells.Item(34, 23).Value = "'Ljustorp"
$ws.Cells.Item(34, 24).Value = "'2025_0528"
$ws.Cells.Item(34, 25).Value = "'2025-06-26"
$ws.Cells.Item(34, 26).Value = "'09:15"
$ws.Cells.Item(34, 27).Value = "'2025-06-26"
$ws.Cells.Item(34, 28).Value = "'09:15"
$ws.Cells.Item(34, 29).Value = "'Brandpåverkad tallstubbe"
$ws.Cells.Item(34, 30).Value = $false
$ws.Cells.Item(34, 31).Value = $false
$ws.Cells.Item(34, 33).Value = $false
$ws.Cells.Item(34, 46).Value = "'"
$ws.Cells.Item(34, 49).Value = "'David Isaksson"
$ws.Cells.Item(34, 50).Value = "'Karin Halldin"
$ws.Cells.Item(34, 51).Value = "'Kustpaketet"

# Row 35
$ws.Cells.Item(35, 1).Value = 131106659
$ws.Cells.Item(35, 2).Value = 57884
$ws.Cells.Item(35, 4).Value = "'NT"
$ws.Cells.Item(35, 5).Value = 100109
$ws.Cells.Item(35, 6).Value = "'Tretåig hackspett"
$ws.Cells.Item(35, 7).Value = "'Picoides tridactylus"
$ws.Cells.Item(35, 8).Value = "'(Linnaeus, 1758)"
$ws.Cells.Item(35, 9).Value = "'"
$ws.Cells.Item(35, 16).Value = "'Lappberget, Mpd"
$ws.Cells.Item(35, 17).Value = 601218
$ws.Cells.Item(35, 18).Value = 6959810
$ws.Cells.Item(35, 19).Value = 10
$ws.Cells.Item(35, 20).Value = "'Västernorrland"
$ws.Cells.Item(35, 21).Value = "'Timrå"
$ws.Cells.Item(35, 22).Value = "'Medelpad"
$ws.Cells.Item(35, 23).Value = "'Ljustorp"
$ws.Cells.Item(35, 24).Value = "'2025_0517"
$ws.Cells.Item(35, 25).Value = "'2025-06-26"
$ws.Cells.Item(35, 26).Value = "'08:33"
$ws.Cells.Item(35, 27).Value = "'2025-06-26"
$ws.Cells.Item(35, 28).Value = "'08:33"
$ws.Cells.Item(35, 29).Value = "'På gran i barrblandskog. Ev. liten hackspett, men är med stor sannolikhet tretåig hackspett enligt diskussion med Anders Forsberg."
$ws.Cells.Item(35, 30).Value = $false
$ws.Cells.Item(35, 31).Value = $false
$ws.Cells.Item(35, 33).Value = $false
$ws.Cells.Item(35, 46).Value = "'"
$ws.Cells.Item(35, 49).Value = "'David Isaksson"
$ws.Cells.Item(35, 50).Value = "'Karin Halldin"
$ws.Cells.Item(35, 51).Value = "'Kustpaketet"

# Row 36
$ws.Cells.Item(36, 1).Value = 131106661
$ws.Cells.Item(36, 2).Value = 79862
$ws.Cells.Item(36, 4).Value = "'NT"
$ws.Cells.Item(36, 5).Value = 6453
$ws.Cells.Item(36, 6).Value = "'Vedskivlav"
$ws.Cells.Item(36, 7).Value = "'Hertelidea botryosa"
$ws.Cells.Item(36, 8).Value = "'(Fr.) Printzen & Kantvilas"
$ws.Cells.Item(36, 9).Value = "'2"
$ws.Cells.Item(36, 10).Value = "'dm²"
$ws.Cells.Item(36, 16).Value = "'Lappberget, Mpd"
$ws.Cells.Item(36, 17).Value = 601240
$ws.Cells.Item(36, 18).Value = 6959782
$ws.Cells.Item(36, 19).Value = 10
$ws.Cells.Item(36, 20).Value = "'Västernorrland"
$ws.Cells.Item(36, 21).Value = "'Timrå"
$ws.Cells.Item(36, 22).Value = "'Medelpad"
$ws.Cells.Item(36, 23).Value = "'Ljustorp"
$ws.Cells.Item(36, 24).Value = "'2025_0515"
$ws.Cells.Item(36, 25).Value = "'2025-06-26"
$ws.Cells.Item(36, 26).Value = "'08:26"
$ws.Cells.Item(36, 27).Value = "'2025-06-26"
$ws.Cells.Item(36, 28).Value = "'08:26"
$ws.Cells.Item(36, 29).Value = "'tallstubbe"
$ws.Cells.Item(36, 30).Value = $false
$ws.Cells.Item(36, 31).Value = $false
$ws.Cells.Item(36, 33).Value = $false
$ws.Cells.Item(36, 46).Value = "'"
$ws.Cells.Item(36, 49).Value = "'David Isaksson"
$ws.Cells.Item(36, 50).Value = "'Måns Svensson"
$ws.Cells.Item(36, 51).Value = "'Kustpaketet"

# Row 37
$ws.Cells.Item(37, 1).Value = 131106646
$ws.Cells.Item(37, 2).Value = 79243
$ws.Cells.Item(37, 4).Value = "'NT"
$ws.Cells.Item(37, 5).Value = 6425
$ws.Cells.Item(37, 6).Value = "'Garnlav"
$ws.Cells.Item(37, 7).Value = "'Alectoria sarmentosa"
$ws.Cells.Item(37, 8).Value = "'(Ach.) Ach."
$ws.Cells.Item(37, 9).Value = "'"
$ws.Cells.Item(37, 16).Value = "'Lappberget, Mpd"
$ws.Cells.Item(37, 17).Value = 601173
$ws.Cells.Item(37, 18).Value = 6959739
$ws.Cells.Item(37, 19).Value = 10
$ws.Cells.Item(37, 20).Value = "'Västernorrland"
$ws.Cells.Item(37, 21).Value = "'Timrå"
$ws.Cells.Item(37, 22).Value = "'Medelpad"
$ws.Cells.Item(37, 23).Value = "'Ljustorp"
$ws.Cells.Item(37, 24).Value = "'2025_0530"
$ws.Cells.Item(37, 25).Value = "'2025-06-26"
$ws.Cells.Item(37, 26).Value = "'09:36"
$ws.Cells.Item(37, 27).Value = "'2025-06-26"
$ws.Cells.Item(37, 28).Value = "'09:36"
$ws.Cells.Item(37, 29).Value = "'tall"
$ws.Cells.Item(37, 30).Value = $false
$ws.Cells.Item(37, 31).Value = $false
$ws.Cells.Item(37, 33).Value = $false
$ws.Cells.Item(37, 46).Value = "'"
$ws.Cells.Item(37, 49).Value = "'David Isaksson"
$ws.Cells.Item(37, 50).Value = "'Måns Svensson"
$ws.Cells.Item(37, 51).Value = "'Kustpaketet"

# Row 38
$ws.Cells.Item(38, 1).Value = 131106643
$ws.Cells.Item(38, 2).Value = 78646
$ws.Cells.Item(38, 4).Value = "'NT"
$ws.Cells.Item(38, 5).Value = 6437
$ws.Cells.Item(38, 6).Value = "'Blanksvart spiklav"
$ws.Cells.Item(38, 7).Value = "'Calicium denigratum"
$ws.Cells.Item(38, 8).Value = "'(Vain.) Tibell"
$ws.Cells.Item(38, 9).Value = "'1"
$ws.Cells.Item(38, 10).Value = "'cm²"
$ws.Cells.Item(38, 16).Value = "'Lappberget, Mpd"
$ws.Cells.Item(38, 17).Value = 601129
$ws.Cells.Item(38, 18).Value = 6959679
$ws.Cells.Item(38, 19).Value = 10
$ws.Cells.Item(38, 20).Value = "'Västernorrland"
$ws.Cells.Item(38, 21).Value = "'Timrå"
$ws.Cells.Item(38, 22).Value = "'Medelpad"
$ws.Cells.Item(38, 23).Value = "'Ljustorp"
$ws.Cells.Item(38, 24).Value = "'2025_0533"
$ws.Cells.Item(38, 25).Value = "'2025-06-26"
$ws.Cells.Item(38, 26).Value = "'09:55"
$ws.Cells.Item(38, 27).Value = "'2025-06-26"
$ws.Cells.Item(38, 28).Value = "'09:55"
$ws.Cells.Item(38, 30).Value = $false
$ws.Cells.Item(38, 31).Value = $false
$ws.Cells.Item(38, 33).Value = $false
$ws.Cells.Item(38, 46).Value = "'"
$ws.Cells.Item(38, 49).Value = "'David Isaksson"
$ws.Cells.Item(38, 50).Value = "'Samuel Koont"
$ws.Cells.Item(38, 51).Value = "'Kustpaketet"

# Row 39
$ws.Cells.Item(39, 1).Value = 131106649
$ws.Cells.Item(39, 2).Value = 91819
$ws.Cells.Item(39, 4).Value = "'LC"
$ws.Cells.Item(39, 5).Value = 1205
$ws.Cells.Item(39, 6).Value = "'Stor aspticka"
$ws.Cells.Item(39, 7).Value = "'Phellinus populicola"
$ws.Cells.Item(39, 8).Value = "'Niemelä"
$ws.Cells.Item(39, 9).Value = "'1"
$ws.Cells.Item(39, 10).Value = "'mycel"
$ws.Cells.Item(39, 16).Value = "'Lappberget, Mpd"
$ws.Cells.Item(39, 17).Value = 601221
$ws.Cells.Item(39, 18).Value = 6959782
$ws.Cells.Item(39, 19).Value = 10
$ws.Cells.Item(39, 20).Value = "'Västernorrland"
$ws.Cells.Item(39, 21).Value = "'Timrå"
$ws.Cells.Item(39, 22).Value = "'Medelpad"
$ws.Cells.Item(39, 23).Value = "'Ljustorp"
$ws.Cells.Item(39, 24).Value = "'2025_0527"
$ws.Cells.Item(39, 25).Value = "'2025-06-26"
$ws.Cells.Item(39, 26).Value = "'09:07"
$ws.Cells.Item(39, 27).Value = "'2025-06-26"
$ws.Cells.Item(39, 28).Value = "'09:07"
$ws.Cells.Item(39, 29).Value = "'aspstubbe"
$ws.Cells.Item(39, 30).Value = $false
$ws.Cells.Item(39, 31).Value = $false
$ws.Cells.Item(39, 33).Value = $false
$ws.Cells.Item(39, 46).Value = "'"
$ws.Cells.Item(39, 49).Value = "'David Isaksson"
$ws.Cells.Item(39, 50).Value = "'Måns Svensson"
$ws.Cells.Item(39, 51).Value = "'Kustpaketet"

# Row 40
$ws.Cells.Item(40, 1).Value = 131106645
$ws.Cells.Item(40, 2).Value = 98930
$ws.Cells.Item(40, 4).Value = "'LC"
$ws.Cells.Item(40, 5).Value = 219790
$ws.Cells.Item(40, 6).Value = "'Fläcknycklar"
$ws.Cells.Item(40, 7).Value = "'Dactylorhiza maculata"
$ws.Cells.Item(40, 8).Value = "'(L.) Soó"
$ws.Cells.Item(40, 9).Value = "'2"
$ws.Cells.Item(40, 10).Value = "'plantor/tuvor"
$ws.Cells.Item(40, 16).Value = "'Lappberget, Mpd"
$ws.Cells.Item(40, 17).Value = 601136
$ws.Cells.Item(40, 18).Value = 6959685
$ws.Cells.Item(40, 19).Value = 10
$ws.Cells.Item(40, 20).Value = "'Västernorrland"
$ws.Cells.Item(40, 21).Value = "'Timrå"
$ws.Cells.Item(40, 22).Value = "'Medelpad"
$ws.Cells.Item(40, 23).Value = "'Ljustorp"
$ws.Cells.Item(40, 24).Value = "'2025_0531"
$ws.Cells.Item(40, 25).Value = "'2025-06-26"
$ws.Cells.Item(40, 26).Value = "'09:51"
$ws.Cells.Item(40, 27).Value = "'2025-06-26"
$ws.Cells.Item(40, 28).Value = "'09:51"
$ws.Cells.Item(40, 30).Value = $false
$ws.Cells.Item(40, 31).Value = $false
$ws.Cells.Item(40, 33).Value = $false
$ws.Cells.Item(40, 46).Value = "'"
$ws.Cells.Item(40, 49).Value = "'David Isaksson"
$ws.Cells.Item(40, 50).Value = "'Samuel Koont"
$ws.Cells.Item(40, 51).Value = "'Kustpaketet"

